$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2005365.4
$ws.Range("J17").Value = 2005365.4
$ws.Range("L17").Value = 6016096.199999999
$ws.Range("N17").Value = -6016432.199999999
$ws.Range("H109").Value = 17205.791
$ws.Range("J109").Value = 17205.791
$ws.Range("L109").Value = 17205.791
$ws.Range("N109").Value = -19979.791
$ws.Range("H112").Value = 1218.5652
$ws.Range("J112").Value = 1271.35
$ws.Range("L112").Value = 3814.05
$ws.Range("N112").Value = -6030.049999999999
$ws.Range("H113").Value = 7889.6553
$ws.Range("I113").Value = 3488.5557
$ws.Range("J113").Value = 9870.15
$ws.Range("K113").Value = 3488.5557
$ws.Range("L113").Value = 9870.15
$ws.Range("M113").Value = -234.5556999999999
$ws.Range("N113").Value = -16378.15
$ws.Range("H114").Value = 37921.75
$ws.Range("J114").Value = 37921.75
$ws.Range("L114").Value = 37921.75
$ws.Range("N114").Value = -46599.75
$ws.Range("H125").Value = 1978018.5
$ws.Range("I125").Value = 50400
$ws.Range("J125").Value = 2941827.8
$ws.Range("K125").Value = 453600
$ws.Range("L125").Value = 26476450.2
$ws.Range("M125").Value = -451140
$ws.Range("N125").Value = -26481370.2
$ws.Range("H138").Value = 2494.28
$ws.Range("I138").Value = 1535.3096
$ws.Range("J138").Value = 3188.7068
$ws.Range("K138").Value = 4605.9288
$ws.Range("L138").Value = 9566.1204
$ws.Range("M138").Value = 534.0712000000003
$ws.Range("N138").Value = -19846.1204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3055.8
$ws.Range("I2").Value = 3055.8
$ws.Range("K2").Value = 3055.8
$ws.Range("M2").Value = -2942.8
$ws.Range("H32").Value = 3631.36
$ws.Range("I32").Value = 3693.484
$ws.Range("J32").Value = 2806
$ws.Range("K32").Value = 3693.484
$ws.Range("L32").Value = 2806
$ws.Range("M32").Value = -3406.484
$ws.Range("N32").Value = -3380
$ws.Range("H110").Value = 125700
$ws.Range("I110").Value = 250450
$ws.Range("J110").Value = 950
$ws.Range("K110").Value = 250450
$ws.Range("L110").Value = 950
$ws.Range("M110").Value = -248405
$ws.Range("N110").Value = -5040
$ws.Range("H113").Value = 30000
$ws.Range("J113").Value = 30000
$ws.Range("L113").Value = 30000
$ws.Range("N113").Value = -38678
$ws.Range("H116").Value = 3055.8
$ws.Range("I116").Value = 3055.8
$ws.Range("K116").Value = 3055.8
$ws.Range("M116").Value = -761.8000000000002
$ws.Range("H119").Value = 24824.5
$ws.Range("J119").Value = 24824.5
$ws.Range("L119").Value = 24824.5
$ws.Range("N119").Value = -34500.5
$ws.Range("H123").Value = 9500
$ws.Range("J123").Value = 9500
$ws.Range("L123").Value = 9500
$ws.Range("N123").Value = -19300
$ws.Range("H134").Value = 51250
$ws.Range("J134").Value = 51250
$ws.Range("L134").Value = 51250
$ws.Range("N134").Value = -61390

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3055.8
$ws.Range("I3").Value = 3055.8
$ws.Range("K3").Value = 3055.8
$ws.Range("M3").Value = -2941.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3514.6875
$ws.Range("I31").Value = 1556.4117
$ws.Range("J31").Value = 5734.067
$ws.Range("K31").Value = 1556.4117
$ws.Range("L31").Value = 5734.067
$ws.Range("M31").Value = -1261.4117
$ws.Range("N31").Value = -6324.067
$ws.Range("H34").Value = 3514.6875
$ws.Range("I34").Value = 1556.4117
$ws.Range("J34").Value = 5734.067
$ws.Range("K34").Value = 1556.4117
$ws.Range("L34").Value = 5734.067
$ws.Range("M34").Value = -1354.4117
$ws.Range("N34").Value = -6138.067
$ws.Range("H51").Value = 29095.227
$ws.Range("J51").Value = 31254.75
$ws.Range("L51").Value = 31254.75
$ws.Range("N51").Value = -32726.75
$ws.Range("H61").Value = 29095.227
$ws.Range("J61").Value = 31254.75
$ws.Range("L61").Value = 31254.75
$ws.Range("N61").Value = -31950.75
$ws.Range("H140").Value = 67338.16
$ws.Range("J140").Value = 67338.16
$ws.Range("L140").Value = 67338.16
$ws.Range("N140").Value = -77698.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 2286.2856
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2286.2856
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 6858.8568
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -7850.8568
$ws.Range("H122").Value = 949.8570999999999
$ws.Range("I122").Value = 523.7143
$ws.Range("J122").Value = 2228.2856
$ws.Range("K122").Value = 4713.428699999999
$ws.Range("L122").Value = 20054.5704
$ws.Range("M122").Value = -2263.428699999999
$ws.Range("N122").Value = -24954.5704
$ws.Range("H131").Value = 843.1
$ws.Range("J131").Value = 870
$ws.Range("L131").Value = 2610
$ws.Range("N131").Value = -12690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 43410.07
$ws.Range("J138").Value = 43410.07
$ws.Range("L138").Value = 43410.07
$ws.Range("N138").Value = -53690.07
$ws.Range("H140").Value = 39375
$ws.Range("J140").Value = 39375
$ws.Range("L140").Value = 39375
$ws.Range("N140").Value = -49735

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 75105.28999999999
$ws.Range("I7").Value = 86789.5
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 86789.5
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -86677.5
$ws.Range("N7").Value = -5224
$ws.Range("H61").Value = 1218
$ws.Range("I61").Value = 827
$ws.Range("K61").Value = 827
$ws.Range("M61").Value = -625
$ws.Range("H109").Value = 29633.75
$ws.Range("J109").Value = 29633.75
$ws.Range("L109").Value = 29633.75
$ws.Range("N109").Value = -32407.75
$ws.Range("H113").Value = 1218
$ws.Range("I113").Value = 827
$ws.Range("K113").Value = 827
$ws.Range("M113").Value = 1343
$ws.Range("H124").Value = 25000
$ws.Range("J124").Value = 25000
$ws.Range("L124").Value = 25000
$ws.Range("N124").Value = -34820
$ws.Range("H126").Value = 75105.28999999999
$ws.Range("I126").Value = 86789.5
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 260368.5
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -257898.5
$ws.Range("N126").Value = -19940
$ws.Range("H133").Value = 42857.715
$ws.Range("J133").Value = 42857.715
$ws.Range("L133").Value = 42857.715
$ws.Range("N133").Value = -47917.715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 700.5
$ws.Range("I107").Value = 700.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2101.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -181.5
$ws.Range("N107").ClearContents()
$ws.Range("H119").Value = 30716.8
$ws.Range("J119").Value = 30716.8
$ws.Range("L119").Value = 30716.8
$ws.Range("N119").Value = -40392.8
$ws.Range("H130").Value = 21045.8
$ws.Range("J130").Value = 21045.8
$ws.Range("L130").Value = 21045.8
$ws.Range("N130").Value = -31085.8
